$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column before column ET (column 150), shifting ET:FX -> EU:FY
$ws.Columns.Item(150).Insert()

# Populate the newly inserted column (now "ET") with the new date header
# and the "-" placeholder used throughout the rest of the column.
$ws.Range("ET1").Value = "23-dec"

for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 150).Value = "-"
}
